$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.102927207946777
$ws.Range("B1").Value = 3.938858032226562
$ws.Range("C1").Value = 2.801111459732056
$ws.Range("D1").Value = 2.228824377059937
$ws.Range("E1").Value = 1.823323845863342
